$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 69, shifting existing rows 69-155 down to 70-156
$ws.Rows.Item(69).Insert()

# Populate the newly inserted row 69 with the new data record
$ws.Range("A69").Value = 4
$ws.Range("B69").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C69").Value = 'Los Lagos'
$ws.Range("D69").Value = 44638
$ws.Range("E69").Value = 10
$ws.Range("F69").Value = 100112009
$ws.Range("G69").Value = 'Acelga'
$ws.Range("H69").Value = 'Sin especificar'
$ws.Range("I69").Value = 'Primera'
$ws.Range("J69").Value = 60
$ws.Range("K69").Value = 10000
$ws.Range("L69").Value = 10000
$ws.Range("M69").Value = 10000
$ws.Range("N69").Value = '$/docena de atados (12 kilos)'
$ws.Range("O69").Value = 'Región de La Araucanía'
$ws.Range("P69").Value = 833
$ws.Range("Q69").Value = 12
$ws.Range("R69").Value = 'Hortaliza'
